$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the HWM 5 homework note from its old location (row 29 / "Mon 01 Dec 2014")
$ws.Range("D29").Value = ""

# Add a new row 31 for the Final Exam (lecture 30 / "Mon 08 Dec 2014"),
# moving the HWM 5 due date to this new row
$ws.Range("A31").Value = 30
$ws.Range("C31").Value = "Final Exam"
$ws.Range("B31").Value = "Mon 08 Dec 2014"
$ws.Range("D31").Value = "HWM 5"

# Match D31's formatting to the rest of the "Homework" column (copy row 29's style)
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null

# Update the active selection to reflect where editing finished
$ws.Range("A32").Select() | Out-Null
